# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps recorded for the handback status report, for both the
# zh-cn and de-de localization sheets.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 2 and 5 both refer to the 0f9d247f... handoff/handback pair
$wsZh.Range("E2").Value = "2016-03-17 20:16:06"
$wsZh.Range("E5").Value = "2016-03-17 20:16:06"
$wsZh.Range("H2").Value = "2016-03-17 20:16:24"
$wsZh.Range("H5").Value = "2016-03-17 20:16:24"

# de-de sheet: rows 2 and 5 both refer to the 0f9d247f... handoff/handback pair
$wsDe.Range("E2").Value = "2016-03-17 20:16:10"
$wsDe.Range("E5").Value = "2016-03-17 20:16:10"
$wsDe.Range("H2").Value = "2016-03-17 20:16:30"
$wsDe.Range("H5").Value = "2016-03-17 20:16:30"
